# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" row into the "总计" (totals) summary sheet,
#    just above the existing "2022-Q2" row.
# 2. Insert a brand new "2022-Q3" worksheet (fund holdings detail),
#    positioned right after "总计" and before "2022-Q2".

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $val) {
    # Force the cell to be stored as text even when the value looks
    # numeric (fund codes like "012850", formatted numbers like "5.88").
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert the 2022-Q3 summary row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push the existing data rows (currently rows 2-5) down by one.
$summary.Rows.Item(2).Insert()

# Copy formatting (bold/centered/bordered style) of the index column
# from the row that just got shifted down, so the new row matches.
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 2).Style = "Normal"
$summary.Cells.Item(2, 3).Value = 2
$summary.Cells.Item(2, 3).Style = "Normal"
$summary.Cells.Item(2, 4).Value = 0.2
$summary.Cells.Item(2, 4).Style = "Normal"

# Renumber the index column (A) of the rows that shifted down.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------
# Step 2: Add the new "2022-Q3" worksheet, right before "2022-Q2"
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$newws = $wb.Worksheets.Add($q2)
$newws.Name = "2022-Q3"

# Match page margins used by the sibling quarter sheets.
$newws.PageSetup.LeftMargin = 54
$newws.PageSetup.RightMargin = 54
$newws.PageSetup.TopMargin = 72
$newws.PageSetup.BottomMargin = 72
$newws.PageSetup.HeaderMargin = 36
$newws.PageSetup.FooterMargin = 36

# Header row.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newws.Cells.Item(1, $i + 2).Value = $headers[$i]
}
# Pick up the bold/centered header style from the sibling sheet.
$q2.Range("B1:H1").Copy()
$newws.Range("B1:H1").PasteSpecial(-4122)

# Data row 2: 012850 - 中融低碳经济3个月持有期混合A
$newws.Cells.Item(2, 1).Value = 0
$newws.Cells.Item(2, 1).Style = "Normal"
Set-TextCell $newws.Cells.Item(2, 2) "012850"
Set-TextCell $newws.Cells.Item(2, 3) "中融低碳经济3个月持有期混合A"
Set-TextCell $newws.Cells.Item(2, 4) "5.88"
Set-TextCell $newws.Cells.Item(2, 5) "86.46"
Set-TextCell $newws.Cells.Item(2, 6) "2.80"
Set-TextCell $newws.Cells.Item(2, 7) "0.1646"
$newws.Cells.Item(2, 8).Value = 10
$newws.Cells.Item(2, 8).Style = "Normal"

# Data row 3: 012851 - 中融低碳经济3个月持有期混合C
$newws.Cells.Item(3, 1).Value = 1
$newws.Cells.Item(3, 1).Style = "Normal"
Set-TextCell $newws.Cells.Item(3, 2) "012851"
Set-TextCell $newws.Cells.Item(3, 3) "中融低碳经济3个月持有期混合C"
Set-TextCell $newws.Cells.Item(3, 4) "1.23"
Set-TextCell $newws.Cells.Item(3, 5) "86.46"
Set-TextCell $newws.Cells.Item(3, 6) "2.80"
Set-TextCell $newws.Cells.Item(3, 7) "0.0344"
$newws.Cells.Item(3, 8).Value = 10
$newws.Cells.Item(3, 8).Style = "Normal"

# Pick up the bold/centered/bordered index-column style (column A) too.
$q2.Range("A2").Copy()
$newws.Range("A2:A3").PasteSpecial(-4122)
$newws.Cells.Item(2, 1).Value = 0
$newws.Cells.Item(3, 1).Value = 1
